$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.167.56'
$ws.Range('E2').Value = '  +0.89%  '
$ws.Range('D3').Value = '1.802.30'
$ws.Range('E3').Value = '  +2.57%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '336.61'
$ws.Range('E5').Value = '  -0.01%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.000'
$ws.Range('E6').Value = '  -0.05%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4625'
$ws.Range('E7').Value = '  +20.09%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3696'
$ws.Range('E8').Value = '  +9.18%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.22'
$ws.Range('E9').Value = '  +1.03%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.148'
$ws.Range('E10').Value = '  +3.49%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07607'
$ws.Range('E11').Value = '  +5.62%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.002'
$ws.Range('E12').Value = '  -0.15%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.37'
$ws.Range('E13').Value = '  +0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.346'
$ws.Range('E14').Value = '  +3.48%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.451'
$ws.Range('E15').Value = '  +4.28%  '
$ws.Range('D16').Value = '1.800.37'
$ws.Range('E16').Value = '  +2.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001097'
$ws.Range('E17').Value = '  +4.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06716'
$ws.Range('E18').Value = '  +1.64%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '81.91'
$ws.Range('E19').Value = '  +3.47%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '1.000'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '17.51'
$ws.Range('E21').Value = '  +5.78%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.416'
$ws.Range('E22').Value = '  +3.32%  '
$ws.Range('D23').Value = '28.174.34'
$ws.Range('E23').Value = '  +0.84%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.90'
$ws.Range('E24').Value = '  +2.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.411'
$ws.Range('E25').Value = '  +0.98%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '20.79'
$ws.Range('E26').Value = '  +5.19%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '153.21'
$ws.Range('E27').Value = '  +0.62%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.380'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').Value = '2.006.44'
$ws.Range('E29').Value = '  +2.43%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '132.98'
$ws.Range('E30').Value = '  +1.01%  '
$ws.Range('E31').Value = '  -1.36%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.031'
$ws.Range('E32').Value = '  +0.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09538'
$ws.Range('E33').Value = '  +8.13%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.859'
$ws.Range('E34').Value = '  +0.85%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.2229'
$ws.Range('E35').Value = '  +5.92%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06380'
$ws.Range('E36').Value = '  +3.21%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02353'
$ws.Range('E37').Value = '  +2.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '12.07'
$ws.Range('E38').Value = '  -0.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.252'
$ws.Range('E39').Value = '  +2.50%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6648'
$ws.Range('E40').Value = '  +0.90%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.513'
$ws.Range('E41').Value = '  +0.36%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.241'
$ws.Range('E42').Value = '  +2.95%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.098'
$ws.Range('E43').Value = '  +1.74%  '
$ws.Range('B44').Value = 'EnergySwap'
$ws.Range('C44').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '14.11'
$ws.Range('E44').Value = '  +2.67%  '
$ws.Range('B45').Value = 'Frax'
$ws.Range('C45').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.9995'
$ws.Range('E45').Value = '  -0.06%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6110'
$ws.Range('E46').Value = '  +1.64%  '
$ws.Range('E47').Value = '  +0.24%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '130.31'
$ws.Range('E48').Value = '  +3.38%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.055'
$ws.Range('E49').Value = '  +2.65%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.07147'
$ws.Range('E50').Value = '  +2.24%  '
$ws.Range('E51').Value = '  +0.47%  '
